# Form the consolidated report: fix the "Absent" (column H) values so that
# they correctly reflect attendance (Absent = 1 - Real) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
